$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 52: replace the old "Test" placeholder row with the Google DNS alarm entry.
$ws.Range("A52").Value = "Google DNS"
$ws.Range("B52").ClearContents()
$ws.Range("C52").Value = "8.8.4.4"
$ws.Range("D52").ClearContents()
$ws.Range("E52").Value = "8.8.8.8"
$ws.Range("F52").ClearContents()

# Row 53: replace the old "Delhi" placeholder row with the Cloudflare DNS entry.
$ws.Range("A53").Value = "Cloudflare DNS"
$ws.Range("B53").ClearContents()
$ws.Range("C53").Value = "1.1.1.1"
$ws.Range("E53").Value = "1.0.0.1"

# Row 54: replace the old "Banglore" placeholder row with the Quad9 entry.
$ws.Range("A54").Value = "Quad9"
$ws.Range("B54").ClearContents()
$ws.Range("C54").Value = "9.9.9.9"
$ws.Range("E54").Value = "149.112.112.112"

# Row 55 (new): OpenDNS entry.
$ws.Range("A55").Value = "OpenDNS"
$ws.Range("C55").Value = "208.67.222.222"
$ws.Range("E55").Value = "208.67.220.220"

# Row 56 (new): Comodo DNS entry.
$ws.Range("A56").Value = "Comodo DNS"
$ws.Range("C56").Value = "8.26.56.26"
$ws.Range("E56").Value = "8.20.247.20"

# Update the visible selection to match the latest edited cell.
[void]$ws.Range("F54").Select()
